$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Image filename updates (jpg -> png/jpeg renames picked up by the sync)
$ws.Range("G16:G18").Value = "img/4.png"
$ws.Range("G30:G35").Value = "img/7.jpeg"
$ws.Range("G67:G69").Value = "img/15.png"

# Restore the scroll position / selection left by the editing session
$ws.Application.Goto($ws.Range("A38"), $true)
$ws.Rows.Item(10).Select()
